$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the reviewer comment ("Should this be ALL OF ?") left by
#    Liz Stebbins. Comments(1).Delete() removes the comment itself as
#    well as the commentRangeStart/commentRangeEnd/commentReference
#    markers that anchor it in the body text, leaving the anchored
#    text itself untouched.
# ------------------------------------------------------------------
if ($d.Comments.Count -gt 0) {
    $d.Comments(1).Delete()
}

# ------------------------------------------------------------------
# 2. Rework the "any of cv_fall + cv_late_fall, cv_spring, cv_winter"
#    clause into "(cv_fall + cv_late_fall), cv_spring, AND cv_winter"
#    (clarifying that ALL three conditions are required, not any one
#    of them) in the bullet that designates UNKNOWN.
# ------------------------------------------------------------------

# "any of " -> "("
$rng = $d.Content
$rng.Find.Execute("any of cv_fall", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($rng.Find.Found) {
    $repl = $d.Range($rng.Start, $rng.Start + 7)
    $repl.Text = "("
}

# "cv_late_fall, cv_spring, cv_winter are" -> "cv_late_fall), cv_spring, AND cv_winter are"
$rng2 = $d.Content
$rng2.Find.Execute("cv_late_fall, cv_spring, cv_winter are", $true, $false, $false, $false, $false, $true, 1, $false, "cv_late_fall), cv_spring, AND cv_winter are", 2) | Out-Null

# Bold just the newly-inserted "AND"
$rng3 = $d.Content
$rng3.Find.Execute("cv_spring, AND cv_winter", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if ($rng3.Find.Found) {
    # "cv_spring, " is 11 characters long, so "AND" starts right after it
    $andStart = $rng3.Start + 11
    $andRng = $d.Range($andStart, $andStart + 3)
    Write-Host "AND range text: [$($andRng.Text)]"
    $andRng.Bold = 1
}

$d.Comments.Count | Out-Null
Write-Host "Edit complete"
